$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update FLASCO's four Q&A answers (row 2) from "No" to "Yes" versions
$ws.Range("C2").Value = "Yes, FLASCO does encompasses community sites, FLASCO's membership primarily consists of academic and hospital-based oncologists and hematologists."
$ws.Range("D2").Value = "Yes, FLASCO is a professional organization focused on clinical oncology and does play a significant role in shaping state or local policy."
$ws.Range("F2").Value = "Yes, FLASCO does provides support for clinical trial recruitment. FLASCO is an organization that focuses on education, advocacy, and networking for oncology professionals in Florida."
$ws.Range("J2").Value = "Yes, The FLASCO board does includes top therapeutic area experts, most of them are practicing oncologists from various specialties."

# Update membership counts in column B
$ws.Range("B2").Value = 4100
$ws.Range("B3").Value = 3100
$ws.Range("B4").Value = 1200
$ws.Range("B5").Value = 1760
$ws.Range("B6").Value = 4000

# Update the active selection
$ws.Range("C10").Select()
